$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that advances by one day
# (2023-10-03 -> 2023-10-04, serial 45202 -> 45203) for every data row.
$ws.Range("C2:C223").Value = 45203
